# Apply the "model and template with unit, description and enum" edit:
#  - Swap displayed text of H1/I1 ("Result" <-> "SamplePortion")
#  - Append a unit annotation to the #float type markers in H2/I2
#  - Add a new row 3 with French field descriptions for each column A-J

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: swap Result / SamplePortion headers
$ws.Range("H1").Value = "SamplePortion"
$ws.Range("I1").Value = "Result"

# Row 2: annotate the float type with its unit (mg)
$ws.Range("H2").Value = "#float,  unit:mg"
$ws.Range("I2").Value = "#float,  unit:mg"

# Row 3 (new): French descriptions / enum hints for each field
$ws.Range("A3").Value = "#Manipulateur"
$ws.Range("B3").Value = "#Desc:IdentifiantEchantillon"
$ws.Range("C3").Value = "#Date"
$ws.Range("D3").Value = "#ModeOderatoireLaboratoire"
$ws.Range("E3").Value = "#AppareilLogicielCritique"
$ws.Range("F3").Value = "#ProduitCritique"
$ws.Range("G3").Value = "#LieuStockageDonneesBrutes"
$ws.Range("H3").Value = "#PriseEssai"
$ws.Range("I3").Value = "#Resultat"
$ws.Range("J3").Value = "#NuméroLotReactif"
